$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert E6 from the numeric 5 back to the text "five"
$ws.Range("E6").Value = "five"

# Update the active selection to E6 (matches the saved view state)
$ws.Range("E6").Select()
